$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "Extension (code)" row above the old row 8 (Subject/Person) ---
$ws.Rows.Item(8).Insert()

# --- Populate the new row 8 ---
$ws.Range("A8").Value = "Extension (code)"
$ws.Range("B8").Value = "Firearm Purchase Prohibition Code"

$xpath = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/j:CourtOrder[@structures:id=../j:ActivityCourtOrderAssociation/j:CourtOrder/@structures:ref]/me-fpp-codes:FirearmPurchaseProhibitionCode"
$ws.Range("D8").Value = $xpath

# Match the wrap/left-top alignment used throughout the sheet for these three cells.
$ws.Range("A8").WrapText = $true
$ws.Range("A8").HorizontalAlignment = -4131
$ws.Range("A8").VerticalAlignment = -4160

$ws.Range("B8").WrapText = $true
$ws.Range("B8").HorizontalAlignment = -4131
$ws.Range("B8").VerticalAlignment = -4160

$ws.Range("D8").WrapText = $true
$ws.Range("D8").HorizontalAlignment = -4131
$ws.Range("D8").VerticalAlignment = -4160

# The second half of the D8 text (the code element itself) is emphasized with its own run.
$splitAt = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/j:CourtOrder[@structures:id=../j:ActivityCourtOrderAssociation/j:CourtOrder/@structures:ref]".Length
$runLen = "/me-fpp-codes:FirearmPurchaseProhibitionCode".Length
$run = $ws.Range("D8").Characters($splitAt + 1, $runLen)
$run.Font.Name = "Calibri"
$run.Font.Size = 12

# Row 8 is a bit taller than a default row, matching the adjacent mapping rows.
$ws.Rows.Item(8).RowHeight = 56

# --- Column B got a bit wider to fit the new "Firearm Purchase Prohibition Code" label ---
$ws.Columns.Item(2).ColumnWidth = 31.67

# --- The active selection ended up on C8 after the edit ---
$ws.Range("C8").Select()
